{"js": "// 1) \"Curso (semestre ideal): EA (8)\" -> \"Curso (semestre ideal): EA (7)\"\nconst matches = context.document.body.search(\"Curso (semestre ideal): EA (8)\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(\"Curso (semestre ideal): EA (7)\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Remove the \"Requisitos\" heading paragraph and its following\n//    \"LOB1202 ... / LOB1232 ...\" bullet-list paragraph entirely.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  const style = para.style;\n  if (style === \"Heading 2\" && text === \"Requisitos\") {\n    toDelete.push(para);\n  } else if (style === \"List Bullet\" && text.indexOf(\"LOB1202\") !== -1 && text.indexOf(\"LOB1232\") !== -1) {\n    toDelete.push(para);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Curso (semestre ideal): EA (8)\" -> \"Curso (semestre ideal): EA (7)\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Curso (semestre ideal): EA (8)\"\n$find.Replacement.Text = \"Curso (semestre ideal): EA (7)\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Remove the \"Requisitos\" heading paragraph and its following\n#    \"LOB1202 ... / LOB1232 ...\" bullet-list paragraph entirely.\n#    Walk backwards so deleting a paragraph doesn't shift the index of\n#    paragraphs still to be examined.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Range.Style.NameLocal\n    $text = $p.Range.Text\n\n    if ($styleName -eq \"Heading 2\" -and $text.Trim() -eq \"Requisitos\") {\n        $p.Range.Delete()\n    }\n    elseif ($styleName -eq \"List Bullet\" -and $text -like \"*LOB1202*\" -and $text -like \"*LOB1232*\") {\n        $p.Range.Delete()\n    }\n}\n"}
